$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cell K1 with the same style (bold/centered/bordered) as the other headers
$ws.Range("A1").Copy($ws.Range("K1"))
$ws.Range("K1").Value = "intervention_type"

# Populate the new intervention_type column for each data row
$ws.Range("K2").Value = "OTHER"
$ws.Range("K3").Value = "DRUG"
$ws.Range("K4").Value = "PROCEDURE"
$ws.Range("K5").Value = "OTHER"
$ws.Range("K6").Value = "PROCEDURE"
$ws.Range("K7").Value = "PROCEDURE"
$ws.Range("K8").Value = "DRUG"
$ws.Range("K9").Value = "PROCEDURE"
$ws.Range("K10").Value = "PROCEDURE"
$ws.Range("K11").Value = "DEVICE"
$ws.Range("K12").Value = "BIOLOGICAL"
$ws.Range("K13").Value = "DRUG"
$ws.Range("K14").Value = "DEVICE"
$ws.Range("K15").Value = "DEVICE"
$ws.Range("C16").Copy($ws.Range("K16"))  # row 16: leave intervention_type empty (matches blank source cells in this row)
$ws.Range("K17").Value = "PROCEDURE"
$ws.Range("K18").Value = "DRUG"
$ws.Range("K19").Value = "RADIATION"
$ws.Range("K20").Value = "BIOLOGICAL"
$ws.Range("K21").Value = "OTHER"
$ws.Range("K22").Value = "DRUG"
$ws.Range("K23").Value = "DRUG"
$ws.Range("C24").Copy($ws.Range("K24"))  # row 24: leave intervention_type empty (matches blank source cells in this row)
$ws.Range("K25").Value = "BIOLOGICAL"
$ws.Range("K26").Value = "OTHER"
$ws.Range("K27").Value = "PROCEDURE"
$ws.Range("K28").Value = "DEVICE"
$ws.Range("K29").Value = "PROCEDURE"
$ws.Range("K30").Value = "DEVICE"
$ws.Range("K31").Value = "OTHER"
$ws.Range("K32").Value = "PROCEDURE"
$ws.Range("K33").Value = "DRUG"
$ws.Range("K34").Value = "BIOLOGICAL"
$ws.Range("K35").Value = "BIOLOGICAL"
$ws.Range("K36").Value = "BIOLOGICAL"

Write-Output "done"
